$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("solar")

$values = @(
    "CHE_11",
    "CHE_15",
    "CHE_25",
    "CHE_3",
    "CHE_21",
    "CHE_9",
    "CHE_4",
    "CHE_7",
    "CHE_14",
    "CHE_18",
    "CHE_13",
    "CHE_24",
    "CHE_5",
    "CHE_8",
    "CHE_10",
    "CHE_22",
    "CHE_20",
    "CHE_1",
    "CHE_6",
    "CHE_0",
    "CHE_17",
    "CHE_19",
    "CHE_12"
)

$startRow = 4
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("AG$row").Value = $values[$i]
}
